$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values for rows 2-12 (columns A-D); column E (dataset name) unchanged.
$data = @{
    2  = @(109.4167172908783, 170.09396004676819, 14.934543132781981, 170.93715929985049)
    3  = @(1.356091976165771, 24.786891937255859, 1.033824682235718, 3.1565861701965332)
    4  = @(26.882259368896481, 51.845715522766113, 11.39080333709717, 21.514628648757931)
    5  = @(11.24188494682312, 86.12848949432373, 6.058922290802002, 18.243011236190799)
    6  = @(3.5003128051757808, 29.67258620262146, 0.89464473724365234, 104.5816335678101)
    7  = @(11.491165637969971, 72.039641618728638, 6.0994770526885986, 16.385441780090328)
    8  = @(431.60133337974548, 108.61535978317259, 7.3583984375, 229.35275721549991)
    9  = @(63.08512544631958, 146.27328681945801, 14.41402983665466, 226.2061040401459)
    10 = @(19.56919693946838, 90.371308565139771, 3.329658985137939, 67.287684679031372)
    11 = @(61.801605224609382, 124.2164630889893, 18.72080397605896, 212.26496911048889)
    12 = @(552.59260940551758, 639.35310506820679, 226.58348870277399, 5253.7760457992554)
    13 = @(127.3657405376434, 173.6004521846771, 124.23669981956481, 2793.7186591625209)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

# Row 13 previously had a distinct style (s="2") applied to all its cells (A13:E13).
# Clear that formatting so the row matches the plain default style like the other rows.
$ws.Range("A13:E13").ClearFormats()
